$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 9308
$ws.Range("F6").Value = 283
$ws.Range("F9").Value = 683
$ws.Range("F13").Value = 308
$ws.Range("F15").Value = 62
$ws.Range("F16").Value = 1571
$ws.Range("F17").Value = 1350
$ws.Range("F19").Value = 57
$ws.Range("F20").Value = 1431
$ws.Range("F22").Value = 269
$ws.Range("F24").Value = 109
$ws.Range("F27").Value = 340
$ws.Range("F28").Value = 340
$ws.Range("F29").Value = 1091
$ws.Range("F30").Value = 14
$ws.Range("F32").Value = 252
$ws.Range("F33").Value = 235
$ws.Range("F34").Value = 67
$ws.Range("F38").Value = 148
$ws.Range("F41").Value = 149
$ws.Range("F42").Value = 555
$ws.Range("F43").Value = 1244
$ws.Range("F44").Value = 710
$ws.Range("F45").Value = 254

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 163
$ws.Range("F13").Value = 42
$ws.Range("F16").Value = 12
$ws.Range("F19").Value = 952
$ws.Range("F20").Value = 27
$ws.Range("F21").Value = 1046
$ws.Range("F22").Value = 257
$ws.Range("F23").Value = 647
$ws.Range("F25").Value = 272
$ws.Range("F26").Value = 272

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 345
$ws.Range("F7").Value = 2288
$ws.Range("F8").Value = 3405

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 9308
$ws.Range("F7").Value = 345
$ws.Range("F8").Value = 3405
$ws.Range("F9").Value = 683
$ws.Range("F13").Value = 1571
$ws.Range("F15").Value = 1350
$ws.Range("F17").Value = 57
$ws.Range("F18").Value = 1431
$ws.Range("F20").Value = 269
$ws.Range("F21").Value = 109
$ws.Range("F23").Value = 340
$ws.Range("F24").Value = 14
$ws.Range("F25").Value = 12
$ws.Range("F29").Value = 952
$ws.Range("F30").Value = 27
$ws.Range("F31").Value = 1046
$ws.Range("F32").Value = 257
$ws.Range("F35").Value = 148
$ws.Range("F37").Value = 272
$ws.Range("F39").Value = 555
$ws.Range("F40").Value = 710

